$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-YearText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
$ws.Range("B2").Value = "On 3D simultaneous attack against manoeuvring target with communication delays"
$ws.Range("C2").Value = "Zhaohui Liu, Yuezu Lv, Jialing Zhou, Liang Hu"
Set-YearText "D2" "2020"
$ws.Range("E2").Value = "10.1177/1729881419894808"
$ws.Range("F2").Value = "Open Access"

# Row 3
$ws.Range("B3").Value = "Ontological security, cyber technology, and states’ responses"
$ws.Range("C3").Value = "Amir Lupovici"
Set-YearText "D3" "2023"
$ws.Range("E3").Value = "10.1177/13540661221130958"
$ws.Range("F3").Value = "Open Access"
$ws.Range("G3").Value = 0

# Row 4
$ws.Range("B4").Value = "A virtual necessity: Some modest steps toward greater cybersecurity"
$ws.Range("C4").Value = "Herbert Lin"
Set-YearText "D4" "2012"
$ws.Range("E4").Value = "10.1177/0096340212459039"
$ws.Range("F4").Value = "Restricted"
$ws.Range("G4").Value = 0

# Row 5
$ws.Range("B5").Value = "A Multi-Phase Network Situational Awareness Cognitive Task Analysis"
$ws.Range("C5").Value = "Robert F. Erbacher, Deborah A. Frincke, Pak Chung Wong, Sarah Moody, Glenn Fink"
Set-YearText "D5" "2010"
$ws.Range("E5").Value = "10.1057/ivs.2010.5"
$ws.Range("G5").Value = 0

# Row 6
$ws.Range("B6").Value = "An Adversarial Model for Expressing Attacks on Control Protocols"
$ws.Range("C6").Value = "Jonathan Butts, Mason Rice, Sujeet Shenoi"
Set-YearText "D6" "2012"
$ws.Range("E6").Value = "10.1177/1548512911449409"

# Row 7
$ws.Range("B7").Value = "The dynamics of cyber conflict between rival antagonists, 2001–11"
$ws.Range("C7").Value = "Brandon Valeriano, Ryan C Maness"
Set-YearText "D7" "2014"
$ws.Range("E7").Value = "10.1177/0022343313518940"
$ws.Range("G7").Value = 1

# Row 8
$ws.Range("B8").Value = "On domains: Cyber and the practice of warfare"
$ws.Range("C8").Value = "Chris McGuffin, Paul Mitchell"
Set-YearText "D8" "2014"
$ws.Range("E8").Value = "10.1177/0020702014540618"
$ws.Range("G8").Value = 1

# Row 9
$ws.Range("B9").Value = "Attrition rates and maneuver in agent-based simulation models"
$ws.Range("C9").Value = "David Ormrod, Benjamin Turnbull"
Set-YearText "D9" "2017"
$ws.Range("E9").Value = "10.1177/1548512917692693"
$ws.Range("G9").Value = 0

# Row 10
$ws.Range("B10").Value = "Wargaming the use of intermediate force capabilities in the gray zone"
$ws.Range("C10").Value = "Kyle D Christensen, Peter Dobias"
Set-YearText "D10" "2024"
$ws.Range("E10").Value = "10.1177/15485129211010227"
$ws.Range("G10").Value = 2

# Row 11
$ws.Range("B11").Value = "Prioritizing investment in military cyber capability using risk analysis"
$ws.Range("C11").Value = "Cayt Rowe, Hossein Seif Zadeh, Ivan L. Garanovich, Li Jiang, Daniel Bilusich, Rick Nunes-Vaz, Anthony Ween"
Set-YearText "D11" "2019"
$ws.Range("E11").Value = "10.1177/1548512917707077"
$ws.Range("F11").Value = "Restricted"
$ws.Range("G11").Value = 0
